$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (column G) values, rows 2-40, to reflect new strikeout
# counts (K) replacing previous Strike# counts.
$kValues = @{
    2  = 5
    3  = 4
    4  = 1
    5  = 2
    6  = 4
    7  = 5
    8  = 5
    9  = 5
    10 = 9
    11 = 4
    12 = 7
    13 = 6
    14 = 5
    15 = 7
    16 = 5
    17 = 7
    18 = 3
    19 = 4
    20 = 8
    21 = 11
    22 = 6
    23 = 6
    24 = 4
    25 = 6
    26 = 7
    27 = 2
    28 = 5
    29 = 7
    30 = 8
    31 = 10
    32 = 7
    33 = 6
    34 = 1
    35 = 1
    36 = 7
    37 = 4
    38 = 3
    39 = 5
    40 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
